$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (price/volume) to Text format so numeric-looking
# strings like "1.002" or "0.1690" are preserved exactly as text, matching
# the source data (which stores these as inline strings, not numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("B2").Value2 = "Bitcoin"
$ws.Range("C2").Value2 = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value2 = "27.661.74"
$ws.Range("E2").Value2 = "  +0.15%  "

$ws.Range("B3").Value2 = "Ethereum"
$ws.Range("C3").Value2 = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value2 = "1.843.19"
$ws.Range("E3").Value2 = "  -0.23%  "

$ws.Range("B4").Value2 = "TetherUSD"
$ws.Range("C4").Value2 = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value2 = "1.002"
$ws.Range("E4").Value2 = "  +0.10%  "

$ws.Range("B5").Value2 = "BNB"
$ws.Range("C5").Value2 = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value2 = "314.98"
$ws.Range("E5").Value2 = "  +0.81%  "

$ws.Range("B6").Value2 = "USDC"
$ws.Range("C6").Value2 = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value2 = "1.002"
$ws.Range("E6").Value2 = "  +0.11%  "

$ws.Range("B7").Value2 = "XRP"
$ws.Range("C7").Value2 = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value2 = "0.4318"
$ws.Range("E7").Value2 = "  +0.71%  "

$ws.Range("B8").Value2 = "Cardano"
$ws.Range("C8").Value2 = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value2 = "0.3704"
$ws.Range("E8").Value2 = "  +1.94%  "

$ws.Range("B9").Value2 = "OKB"
$ws.Range("C9").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value2 = "45.16"
$ws.Range("E9").Value2 = "  +0.47%  "

$ws.Range("B10").Value2 = "Dogecoin"
$ws.Range("C10").Value2 = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value2 = "0.07333"
$ws.Range("E10").Value2 = "  +0.26%  "

$ws.Range("B11").Value2 = "Polygon"
$ws.Range("C11").Value2 = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value2 = "0.8785"
$ws.Range("E11").Value2 = "  +0.31%  "

$ws.Range("B12").Value2 = "Solana"
$ws.Range("C12").Value2 = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value2 = "21.04"
$ws.Range("E12").Value2 = "  +1.81%  "

$ws.Range("B13").Value2 = "WrappedEther"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value2 = "1.917.46"
$ws.Range("E13").Value2 = "  +3.73%  "

$ws.Range("B14").Value2 = "Polkadot"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value2 = "5.475"
$ws.Range("E14").Value2 = "  +2.76%  "

$ws.Range("B15").Value2 = "Chainlink"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value2 = "6.599"
$ws.Range("E15").Value2 = "  +1.24%  "

$ws.Range("B16").Value2 = "TRON"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value2 = "0.06954"
$ws.Range("E16").Value2 = "  +0.60%  "

$ws.Range("B17").Value2 = "BinanceUSD"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value2 = "1.004"
$ws.Range("E17").Value2 = "  +0.28%  "

$ws.Range("B18").Value2 = "Litecoin"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value2 = "81.11"
$ws.Range("E18").Value2 = "  +1.50%  "

$ws.Range("B19").Value2 = "ShibaInu"
$ws.Range("C19").Value2 = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value2 = "0.000009023"
$ws.Range("E19").Value2 = "  -0.19%  "

$ws.Range("B20").Value2 = "Dai"
$ws.Range("C20").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value2 = "1.001"
$ws.Range("E20").Value2 = "  +0.11%  "

$ws.Range("B21").Value2 = "Avalanche"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value2 = "15.55"
$ws.Range("E21").Value2 = "  +1.19%  "

$ws.Range("B22").Value2 = "WrappedBTC"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value2 = "28.008.32"
$ws.Range("E22").Value2 = "  +1.33%  "

$ws.Range("B23").Value2 = "Uniswap"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value2 = "5.126"
$ws.Range("E23").Value2 = "  +3.45%  "

$ws.Range("B24").Value2 = "Cosmos"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value2 = "10.99"
$ws.Range("E24").Value2 = "  +5.76%  "

$ws.Range("B25").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value2 = "2.154.72"
$ws.Range("E25").Value2 = "  +2.21%  "

$ws.Range("B26").Value2 = "Toncoin"
$ws.Range("C26").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value2 = "1.988"
$ws.Range("E26").Value2 = "  -0.06%  "

$ws.Range("B27").Value2 = "Monero"
$ws.Range("C27").Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value2 = "154.11"
$ws.Range("E27").Value2 = "  -0.61%  "

$ws.Range("B28").Value2 = "EthereumClassic"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value2 = "18.93"
$ws.Range("E28").Value2 = "  +1.09%  "

$ws.Range("B29").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value2 = "5.315"
$ws.Range("E29").Value2 = "  +0.66%  "

$ws.Range("B30").Value2 = "BitcoinCash"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value2 = "115.87"
$ws.Range("E30").Value2 = "  -4.79%  "

$ws.Range("B31").Value2 = "LidoDAOToken"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D31").Value2 = "1.877"
$ws.Range("E31").Value2 = "  +1.74%  "

$ws.Range("B32").Value2 = "Stellar"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value2 = "0.08931"
$ws.Range("E32").Value2 = "  +0.28%  "

$ws.Range("B33").Value2 = "ImmutableX"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value2 = "0.7873"
$ws.Range("E33").Value2 = "  +2.96%  "

$ws.Range("B34").Value2 = "Filecoin"
$ws.Range("C34").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value2 = "4.615"
$ws.Range("E34").Value2 = "  +1.34%  "

$ws.Range("B35").Value2 = "ARBITRUM"
$ws.Range("C35").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value2 = "1.177"
$ws.Range("E35").Value2 = "  +6.67%  "

$ws.Range("B36").Value2 = "HuobiToken"
$ws.Range("C36").Value2 = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value2 = "2.975"
$ws.Range("E36").Value2 = "  -0.76%  "

$ws.Range("B37").Value2 = "Frax"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value2 = "1.002"
$ws.Range("E37").Value2 = "  +0.19%  "

$ws.Range("B38").Value2 = "Hedera"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value2 = "0.05439"
$ws.Range("E38").Value2 = "  +0.48%  "

$ws.Range("B39").Value2 = "TrustWalletToken"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value2 = "1.104"
$ws.Range("E39").Value2 = "  +1.48%  "

$ws.Range("B40").Value2 = "VeChain"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value2 = "0.01962"
$ws.Range("E40").Value2 = "  +1.37%  "

$ws.Range("B41").Value2 = "MXToken"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value2 = "2.847"
$ws.Range("E41").Value2 = "  +1.03%  "

$ws.Range("B42").Value2 = "TheSandbox"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value2 = "0.5176"
$ws.Range("E42").Value2 = "  +1.94%  "

$ws.Range("B43").Value2 = "Algorand"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value2 = "0.1690"
$ws.Range("E43").Value2 = "  +2.17%  "

$ws.Range("B44").Value2 = "FraxShare"
$ws.Range("C44").Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value2 = "6.813"
$ws.Range("E44").Value2 = "  +0.82%  "

$ws.Range("B45").Value2 = "Aptos"
$ws.Range("C45").Value2 = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value2 = "8.645"
$ws.Range("E45").Value2 = "  +3.32%  "

$ws.Range("B46").Value2 = "EnergySwap"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value2 = "10.61"
$ws.Range("E46").Value2 = "  +1.98%  "

$ws.Range("B47").Value2 = "Decentraland"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value2 = "0.4779"
$ws.Range("E47").Value2 = "  +2.08%  "

$ws.Range("B48").Value2 = "Quant"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value2 = "106.76"
$ws.Range("E48").Value2 = "  +1.74%  "

$ws.Range("B49").Value2 = "Cronos"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value2 = "0.06554"
$ws.Range("E49").Value2 = "  +0.05%  "

$ws.Range("B50").Value2 = "PaxDollar"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value2 = "1.001"
$ws.Range("E50").Value2 = "  +0.15%  "

$ws.Range("B51").Value2 = "NEARProtocol"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value2 = "1.663"
$ws.Range("E51").Value2 = "  +2.54%  "
